$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = -20.91869999999999
$ws.Range("A6").Value = -22.83590000000001
$ws.Range("A7").Value = -21.78629999999999
$ws.Range("B7").Value = 4.946799999999999
$ws.Range("B12").Value = 4.618599999999998
$ws.Range("C13").Value = -13.68509999999999
$ws.Range("C14").Value = -14.0691
$ws.Range("B15").Value = 4.853699999999996
$ws.Range("A16").Value = -21.58219999999999
$ws.Range("C16").Value = -11.9936
$ws.Range("C19").Value = -11.97250000000001
$ws.Range("A20").Value = -22.76210000000001
$ws.Range("B20").Value = 4.793199999999997
$ws.Range("B21").Value = 10.4216
$ws.Range("B22").Value = 10.2871
$ws.Range("C22").Value = -12.8541
$ws.Range("B23").Value = 9.126499999999998
$ws.Range("A28").Value = -22.00489999999999
$ws.Range("A29").Value = -21.67109999999999
$ws.Range("A32").Value = -21.20589999999999
$ws.Range("B34").Value = 9.342600000000006
$ws.Range("C36").Value = -12.8379
$ws.Range("A40").Value = -19.63709999999999
$ws.Range("B42").Value = 9.972699999999998
$ws.Range("B43").Value = 5.856000000000003
$ws.Range("B44").Value = 4.790500000000002
$ws.Range("B45").Value = 5.2182
$ws.Range("A46").Value = -22.1819
$ws.Range("B46").Value = 5.182599999999993
$ws.Range("C46").Value = -13.4244
$ws.Range("B50").Value = 4.704499999999997
$ws.Range("C50").Value = -13.39999999999999
$ws.Range("A51").Value = -22.2978
$ws.Range("B51").Value = 5.459299999999998
$ws.Range("A52").Value = -22.06550000000001
$ws.Range("A57").Value = -22.81050000000002
$ws.Range("A59").Value = -22.17360000000001
$ws.Range("A62").Value = -22.14870000000001
$ws.Range("A66").Value = -21.55890000000001
$ws.Range("B66").Value = 4.838299999999997
$ws.Range("B67").Value = 5.0749
$ws.Range("A73").Value = -20.2213
$ws.Range("A74").Value = -21.92469999999998
$ws.Range("B79").Value = 9.661700000000002
$ws.Range("B84").Value = 5.372500000000001
$ws.Range("A92").Value = -21.39970000000002
$ws.Range("B92").Value = 5.105499999999996
$ws.Range("C95").Value = -11.6443
$ws.Range("B97").Value = 5.959599999999998
$ws.Range("C97").Value = -11.6604
$ws.Range("A100").Value = -22.1507
